$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 35 (Keysborough) and row 54 (Springvale/IKEA), shifting cells up.
$ws.Rows.Item(54).Delete()
$ws.Rows.Item(35).Delete()
